$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove the old "DSS" source rows (currently rows 26-27); their content will
# be re-created further down the sheet (rows 32-33) to make room for a new
# "Number of employees / Assets / Turnover" breakdown table.
# ---------------------------------------------------------------------------
$ws.Range("A26:D27").Delete() | Out-Null

# ---------------------------------------------------------------------------
# New table header (row 21): Number of employees | Assets | Turnover
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "Number of employees"
$ws.Range("B21").Style = "title"

$ws.Range("C21").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C21").Style = "title"

$ws.Range("D21").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D21").Style = "title"

# ---------------------------------------------------------------------------
# Row 22: Micro
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "Micro"
$ws.Range("A22").Style = "Normal"

# ---------------------------------------------------------------------------
# Row 23: Small
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "Small"
$ws.Range("A23").Style = "Normal"

# ---------------------------------------------------------------------------
# Row 24: Medium
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "Medium"
$ws.Range("A24").Style = "Normal"

$ws.Range("B24").Value = '''=<200 all sectors'
$ws.Range("B24").Style = "Normal"

$ws.Range("D24").Value = '<S$100,000,000'
$ws.Range("D24").Style = "Normal"

# ---------------------------------------------------------------------------
# Row 25: Large
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "Large"
$ws.Range("A25").Style = "Normal"

$ws.Range("B25").Value = '>200'
$ws.Range("B25").Style = "Normal"

$ws.Range("D25").Value = '> S$ 100,000,000'
$ws.Range("D25").Style = "Normal"

# ---------------------------------------------------------------------------
# Re-create the source rows further down the sheet (rows 32-33)
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = "DSS"
$ws.Range("A32").Style = "title"

$ws.Range("A33").Value = "Department of Statistics Singapore.  Subject: Companies and Businesses, Topic: Enterprises, Title: M600981 - Topline Estimates For All Enterprises And SMEs, Annual.  Utilizing SingStat Table Builder. Singapore."
$ws.Range("A33").Style = "source"
